$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that currently have empty placeholder cells in E:G which should be
# removed entirely (no more <c> elements for those columns on these rows).
$rowsToClear = @(2,3,4,5,6,7,8,12,14,15,16)
foreach ($r in $rowsToClear) {
    $rangeAddr = "E${r}:G${r}"
    $ws.Range($rangeAddr).ClearContents()
}

# Row 9 - Business Segment Performance: update description + add examples
$ws.Range("D9").Value = 'Performance by business segment — segment names vary by bank (e.g., Global Banking & Markets, Wholesale Banking, Personal & Commercial). Include revenue, earnings, and growth metrics for each reported division'
$ws.Range("E9").Value = 'BNS: Canadian Banking, International Banking, Global Banking and Markets, Global Wealth Management'
$ws.Range("F9").Value = 'TD: Canadian P&C, U.S. Retail, Wealth & Insurance, Wholesale Banking'
$ws.Range("G9").Value = 'RY: Personal & Commercial, Capital Markets, Wealth Management, Insurance, Investor & Treasury Services'

# Row 10 - Canadian Economic Outlook: update description + add examples
$ws.Range("D10").Value = 'Commentary on Canadian economic environment, Bank of Canada rate outlook, provincial economic trends, and banking industry outlook including competitive dynamics and sector-wide earnings expectations'
$ws.Range("E10").Value = 'Bank of Canada rate trajectory and monetary policy impact'
$ws.Range("F10").Value = 'GDP growth, employment trends, consumer confidence'
$ws.Range("G10").Value = 'Industry-wide credit trends or competitive dynamics'

# Row 11 - Forward Guidance & Outlook: update description + add examples
$ws.Range("D11").Value = 'Management guidance for future quarters, earnings outlook, strategic priorities, growth targets, medium-term objectives, and forward-looking statements'
$ws.Range("E11").Value = 'Net income growth targets or ROE guidance'
$ws.Range("F11").Value = 'Operating leverage and efficiency ratio targets'
$ws.Range("G11").Value = 'Medium-term strategic objectives and capital deployment plans'

# Row 13 - Digital & Technology Initiatives: update description + add examples
$ws.Range("D13").Value = 'Digital transformation efforts, technology investments, AI/ML initiatives, fintech partnerships, and technology-driven productivity and operational efficiency improvements'
$ws.Range("E13").Value = 'Digital adoption rates and online banking migration'
$ws.Range("F13").Value = 'Technology-driven productivity gains and process automation'
$ws.Range("G13").Value = 'Cloud migration, cybersecurity investments, data analytics capabilities'

Write-Host "Edit complete"
